$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 250 (shifts rows 250..359 down to 251..360)
$ws.Rows("250:250").Insert()

# Populate the newly inserted row 250 with the new record's data.
# Columns A, B, C, E, F, G, H, I, R stay consistent with the rest of the
# dataset (same market/region/category/variety/quality/classification).
$ws.Range("A250").Value = 5
$ws.Range("B250").Value = "Macroferia Regional de Talca"
$ws.Range("C250").Value = "Maule"
$ws.Range("D250").Value = 44636
$ws.Range("E250").Value = 7
$ws.Range("F250").Value = 100112043
$ws.Range("G250").Value = "Pepino ensalada"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 350
$ws.Range("K250").Value = 20000
$ws.Range("L250").Value = 20000
$ws.Range("M250").Value = 20000
$ws.Range("N250").Value = "$/caja 80 unidades"
$ws.Range("O250").Value = "Región del Maule"
$ws.Range("P250").Value = 250
$ws.Range("Q250").Value = 80
$ws.Range("R250").Value = "Hortaliza"
